# Atualização de bases das ligas, do dia: 23-02-2024 às 23:34
# Adds the new team "FA Odisha" plus fixture data:
#  - completes row 146 (Hyderabad FC vs Mumbai City FC) with its final score/odds
#  - replaces row 147 with a new upcoming fixture (FA Odisha vs Mohun Bagan SG)
#  - pushes the former row 147 fixture (Bengaluru vs Hyderabad FC) down to row 148
#    with updated odds

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 146: Hyderabad FC vs Mumbai City FC - fill in the result + odds that
# were missing, and update the closing odds columns.
# ---------------------------------------------------------------------------
$ws.Cells.Item(146, 8).Value = 0          # H146 FTHG
$ws.Cells.Item(146, 9).Value = 2          # I146 FTAG
$ws.Cells.Item(146, 10).Value = "A"       # J146 FTR

$ws.Cells.Item(146, 14).Value = 3.8       # N146
$ws.Cells.Item(146, 15).Value = 3.75      # O146
$ws.Cells.Item(146, 16).Value = 1.85      # P146
$ws.Cells.Item(146, 18).Value = 1.975     # R146
$ws.Cells.Item(146, 19).Value = 1.875     # S146
$ws.Cells.Item(146, 22).Value = 2         # V146
$ws.Cells.Item(146, 23).Value = -1        # W146
$ws.Cells.Item(146, 24).Value = -1        # X146
$ws.Cells.Item(146, 25).Value = 0.8500000000000001  # Y146
$ws.Cells.Item(146, 26).Value = -1        # Z146
$ws.Cells.Item(146, 27).Value = 0.875     # AA146
$ws.Cells.Item(146, 28).Value = -1        # AB146
$ws.Cells.Item(146, 29).Value = 1         # AC146

# ---------------------------------------------------------------------------
# Row 148 (new): move the former row-147 fixture (Bengaluru vs Hyderabad FC,
# id 7751751) one row down, with refreshed odds.
# ---------------------------------------------------------------------------
# Row 148 is brand new, so first clone the column formatting used throughout
# the table (bold/bordered id column, date-formatted date column) from the
# row above before writing values into it.
$ws.Cells.Item(147, 1).Copy()
$ws.Cells.Item(148, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(147, 5).Copy()
$ws.Cells.Item(148, 5).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(148, 1).Value = 146
$ws.Cells.Item(148, 2).Value = 7751751
$ws.Cells.Item(148, 3).Value = "India Super League"
$ws.Cells.Item(148, 4).Value = "India Super League"
$ws.Cells.Item(148, 5).Value = 45346.45833333334
$ws.Cells.Item(148, 6).Value = "Bengaluru"
$ws.Cells.Item(148, 7).Value = "Hyderabad FC"

$ws.Cells.Item(148, 11).Value = 1.2       # K148
$ws.Cells.Item(148, 12).Value = 5.5       # L148
$ws.Cells.Item(148, 13).Value = 10        # M148
$ws.Cells.Item(148, 14).Value = 1.142     # N148
$ws.Cells.Item(148, 15).Value = 6.5       # O148
$ws.Cells.Item(148, 16).Value = 17        # P148
$ws.Cells.Item(148, 17).Value = -2        # Q148
$ws.Cells.Item(148, 18).Value = 1.85      # R148
$ws.Cells.Item(148, 19).Value = 1.95      # S148
$ws.Cells.Item(148, 20).Value = 3         # T148
$ws.Cells.Item(148, 21).Value = 1.975     # U148
$ws.Cells.Item(148, 22).Value = 1.825     # V148
$ws.Cells.Item(148, 23).Value = 0         # W148
$ws.Cells.Item(148, 24).Value = 0         # X148
$ws.Cells.Item(148, 25).Value = 0         # Y148
$ws.Cells.Item(148, 26).Value = 0         # Z148
$ws.Cells.Item(148, 27).Value = 0         # AA148

# ---------------------------------------------------------------------------
# Row 147: overwrite with the brand new fixture, FA Odisha vs Mohun Bagan SG
# (id 7873049), replacing the old Bengaluru vs Hyderabad FC data that now
# lives in row 148.
# ---------------------------------------------------------------------------
$ws.Cells.Item(147, 2).Value = 7873049
$ws.Cells.Item(147, 5).Value = 45346.35416666666
$ws.Cells.Item(147, 6).Value = "FA Odisha"
$ws.Cells.Item(147, 7).Value = "Mohun Bagan SG"

$ws.Cells.Item(147, 11).Value = 2.25      # K147
$ws.Cells.Item(147, 12).Value = 3.5       # L147
$ws.Cells.Item(147, 13).Value = 2.625     # M147
$ws.Cells.Item(147, 14).Value = 2.375     # N147
$ws.Cells.Item(147, 15).Value = 3.5       # O147
$ws.Cells.Item(147, 16).Value = 2.5       # P147
$ws.Cells.Item(147, 17).Value = 0         # Q147
$ws.Cells.Item(147, 18).Value = 1.8       # R147
$ws.Cells.Item(147, 19).Value = 2         # S147
$ws.Cells.Item(147, 20).Value = 2.5       # T147
$ws.Cells.Item(147, 21).Value = 1.85      # U147
$ws.Cells.Item(147, 22).Value = 1.95      # V147
$ws.Cells.Item(147, 23).Value = 0         # W147
$ws.Cells.Item(147, 24).Value = 0         # X147
$ws.Cells.Item(147, 25).Value = 0         # Y147
$ws.Cells.Item(147, 26).Value = 0         # Z147
$ws.Cells.Item(147, 27).Value = 0         # AA147
